$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'sichern',
    'saufen',
    'wüten',
    'füllen',
    'laufen',
    'stehlen',
    'schreiten',
    'führen',
    'kehren',
    'ändern',
    'merken',
    'bremsen',
    'altern',
    'erben',
    'rufen',
    'feiern',
    'runden',
    'scheinen',
    'schützen',
    'hauen',
    'klappen',
    'warten',
    'greifen',
    'rauchen',
    'brauchen',
    'schulden',
    'rühren',
    'ahnen',
    'stoppen',
    'lesen',
    'fließen',
    'wählen',
    'reizen',
    'arten',
    'bitten',
    'spüren',
    'wehtun',
    'kümmern',
    'bieten',
    'filmen',
    'klettern',
    'räumen',
    'sprengen',
    'zielen',
    'tauchen',
    'rauben',
    'jubeln',
    'töten',
    'platzen',
    'orten',
    'geben',
    'öffnen',
    'flehen',
    'werfen',
    'rechnen',
    'achten',
    'sitzen',
    'fangen',
    'formen',
    'rücken',
    'weichen',
    'stillen',
    'dienen',
    'meistern',
    'biegen',
    'trauen',
    'regnen',
    'fischen',
    'wiegen',
    'prügeln',
    'landen',
    'suchen',
    'fahren',
    'graben',
    'fällen',
    'ziehen',
    'loben',
    'schnellen',
    'seufzen',
    'heilen',
    'dringen',
    'spielen',
    'decken',
    'wetten',
    'knarren',
    'mögen',
    'schwingen',
    'garen',
    'ehren',
    'betteln',
    'quälen',
    'münzen',
    'äußern',
    'schrecken',
    'hupen',
    'jagen',
    'sterben',
    'krachen',
    'warnen',
    'schwören',
    'wundern',
    'folgen',
    'plaudern',
    'deuten',
    'klingen',
    'schlucken',
    'irren',
    'strahlen',
    'stecken',
    'pfeifen',
    'löschen',
    'drohen',
    'wehren',
    'hindern',
    'bergen',
    'stimmen',
    'planen',
    'kichern',
    'lösen',
    'spinnen',
    'enden',
    'malen',
    'bellen',
    'pflanzen',
    'handeln',
    'lächeln',
    'trennen',
    'hören',
    'liegen',
    'beißen',
    'flüchten',
    'wecken',
    'machen',
    'mauern',
    'gründen',
    'wirken',
    'wahren',
    'wenden',
    'mühen',
    'zünden',
    'lügen',
    'sperren',
    'lassen',
    'narren',
    'schwimmen',
    'zögern',
    'zeigen',
    'ächzen',
    'werden',
    'heulen',
    'siegen',
    'wandern',
    'lenken',
    'wachsen',
    'streichen',
    'kennen',
    'freuen',
    'liefern',
    'leisten',
    'bluten',
    'gelten',
    'tollen',
    'kranken',
    'schwächen',
    'ärgern',
    'sinken',
    'schmecken',
    'kosten',
    'backen',
    'doppeln',
    'drehen',
    'lohnen',
    'tropfen',
    'scheitern',
    'bauen',
    'heben',
    'nerven',
    'schenken',
    'hassen',
    'lockern',
    'tauschen',
    'sorgen',
    'stammen',
    'grüßen',
    'helfen',
    'läuten',
    'pissen',
    'schleppen',
    'fallen',
    'rasen',
    'boxen',
    'treiben'
)

for ($i = 0; $i -lt $words.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
